# Updates cryptos list values (price/volume) and reorders two coin rows,
# reproducing the data refresh captured in the source OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.107.75'
$ws.Range('E2').Value = '  +5.52%  '
# Row 3
$ws.Range('D3').Value = '1.917.28'
$ws.Range('E3').Value = '  +2.34%  '
# Row 4
$ws.Range('E4').Value = '  -0.85%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '329.77'
$ws.Range('E5').Value = '  +4.48%  '
# Row 6
$ws.Range('E6').Value = '  -0.75%  '
# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5208'
$ws.Range('E7').Value = '  +2.31%  '
# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4084'
$ws.Range('E8').Value = '  +4.51%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08548'
$ws.Range('E9').Value = '  +2.18%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.90'
$ws.Range('E10').Value = '  +2.39%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.125'
$ws.Range('E11').Value = '  +1.65%  '
# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.28'
$ws.Range('E12').Value = '  +9.14%  '
# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.414'
$ws.Range('E13').Value = '  +3.28%  '
# Row 14
$ws.Range('D14').Value = '1.932.24'
$ws.Range('E14').Value = '  +3.41%  '
# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.394'
$ws.Range('E15').Value = '  +1.66%  '
# Row 16
$ws.Range('E16').Value = '  -0.94%  '
# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '95.32'
$ws.Range('E17').Value = '  +4.49%  '
# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001115'
$ws.Range('E18').Value = '  +1.29%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06689'
$ws.Range('E19').Value = '  -0.80%  '
# Row 20
$ws.Range('E20').Value = '  +4.28%  '
# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.74%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.008'
$ws.Range('E22').Value = '  +1.43%  '
# Row 23
$ws.Range('D23').Value = '30.118.77'
$ws.Range('E23').Value = '  +5.36%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.36'
$ws.Range('E24').Value = '  +2.37%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.199'
$ws.Range('E25').Value = '  -0.61%  '
# Row 26
$ws.Range('D26').Value = '2.149.49'
$ws.Range('E26').Value = '  +3.22%  '
# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '159.87'
$ws.Range('E27').Value = '  +1.17%  '
# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '21.03'
$ws.Range('E28').Value = '  +2.12%  '
# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.438'
$ws.Range('E29').Value = '  +0.79%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '128.93'
$ws.Range('E30').Value = '  +1.62%  '
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.078'
$ws.Range('E31').Value = '  +3.09%  '
# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1063'
$ws.Range('E32').Value = '  +2.31%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.040'
$ws.Range('E33').Value = '  +5.39%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.636'
$ws.Range('E34').Value = '  +0.48%  '
# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02485'
$ws.Range('E35').Value = '  +0.99%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06605'
$ws.Range('E36').Value = '  +0.26%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2208'
$ws.Range('E37').Value = '  +2.01%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.231'
$ws.Range('E38').Value = '  +4.17%  '
# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.172'
$ws.Range('E39').Value = '  +2.56%  '
# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.918'
$ws.Range('E40').Value = '  -0.09%  '
# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6532'
$ws.Range('E41').Value = '  +2.62%  '
# Row 42
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.65'
$ws.Range('E42').Value = '  +4.92%  '
# Row 43
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.245'
$ws.Range('E43').Value = '  +0.50%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6150'
$ws.Range('E44').Value = '  +2.47%  '
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.31'
$ws.Range('E45').Value = '  +2.17%  '
# Row 46
$ws.Range('E46').Value = '  +2.06%  '
# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.076'
$ws.Range('E47').Value = '  +3.57%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.246'
$ws.Range('E48').Value = '  +2.74%  '
# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '124.41'
$ws.Range('E49').Value = '  +1.47%  '
# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.161'
$ws.Range('E50').Value = '  +9.01%  '
# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '79.74'
$ws.Range('E51').Value = '  +4.37%  '
